$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 10.402079
$ws.Range("N2").Value = 20.804158
$ws.Range("O2").Value = 0.117441350183963
$ws.Range("P2").Value = 0.08862141909929068
$ws.Range("Q2").Value = 0.04605347109266667
$ws.Range("R2").Value = 0.276320826556
$ws.Range("S2").Value = 0.117441350183963
$ws.Range("T2").Value = 0.08862141909929068

# Row 3
$ws.Range("O3").Value = 0.6341369869521791
$ws.Range("P3").Value = 0.7177810830557603
$ws.Range("Q3").Value = 0.2486705862257778
$ws.Range("R3").Value = 2.238035276032
$ws.Range("S3").Value = 0.6341369869521791
$ws.Range("T3").Value = 0.7177810830557603

# Row 4
$ws.Range("M4").Value = 0.2909543333333333
$ws.Range("N4").Value = 0.8728629999999999
$ws.Range("O4").Value = 0.003284926960133785
$ws.Range("P4").Value = 0.003718216220971988
$ws.Range("Q4").Value = 0.001288151818444444
$ws.Range("R4").Value = 0.011593366366
$ws.Range("S4").Value = 0.003284926960133785
$ws.Range("T4").Value = 0.003718216220971988

# Row 5
$ws.Range("M5").Value = 20.562391
$ws.Range("N5").Value = 41.124782
$ws.Range("O5").Value = 0.232153107282743
$ws.Range("P5").Value = 0.175183083160057
$ws.Range("Q5").Value = 0.09103655908733332
$ws.Range("R5").Value = 0.5462193545239999
$ws.Range("S5").Value = 0.232153107282743
$ws.Range("T5").Value = 0.175183083160057

# Row 6
$ws.Range("M6").Value = 0.5741476666666667
$ws.Range("N6").Value = 1.722443
$ws.Range("O6").Value = 0.006482230828885768
$ws.Range("P6").Value = 0.007337251667557973
$ws.Range("Q6").Value = 0.002541943102888889
$ws.Range("R6").Value = 0.022877487926
$ws.Range("S6").Value = 0.006482230828885768
$ws.Range("T6").Value = 0.007337251667557973

# Row 7
$ws.Range("M7").Value = 0.5758453333333333
$ws.Range("N7").Value = 1.727536
$ws.Range("O7").Value = 0.006501397792095299
$ws.Range("P7").Value = 0.00735894679636216
$ws.Range("Q7").Value = 0.002549459239111111
$ws.Range("R7").Value = 0.022945133152
$ws.Range("S7").Value = 0.006501397792095299
$ws.Range("T7").Value = 0.00735894679636216
